# Update "paises" workbook (COVID-19 country stats) with new data pull.
# - Refresh the "Datos actualizados..." timestamp in A1
# - Update numeric stats for countries whose figures changed
# - Re-rank country rows whose totals (column B) now overtake their
#   neighbour, swapping the two rows (name + full data row move together)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp update (A1) ---------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 4 de Julio de 2020 a las 13:47"

# --- Plain value refreshes (no re-ranking needed) -----------------------

# Row 4: Estados Unidos
$ws.Range("B4").Value = 2891267
$ws.Range("C4").Value = 679
$ws.Range("D4").Value = 1235965
$ws.Range("E4").Value = 1523190
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 11
$ws.Range("H4").Value = 132112

# Row 5: Brasil
$ws.Range("B5").Value = 1545458
$ws.Range("C5").Value = 2117
$ws.Range("D5").Value = 978615
$ws.Range("E5").Value = 503548
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 41
$ws.Range("H5").Value = 63295

# Row 21: Banglades
$ws.Range("B21").Value = 159679
$ws.Range("C21").Value = 3288
$ws.Range("D21").Value = 70721
$ws.Range("E21").Value = 86961
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 29
$ws.Range("H21").Value = 1997

# Row 24: Catar
$ws.Range("B24").Value = 99183
$ws.Range("C24").Value = 530
$ws.Range("D24").Value = 90387
$ws.Range("E24").Value = 8673
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 2
$ws.Range("H24").Value = 123

# Row 96: Hungria
$ws.Range("D96").Value = 2784
$ws.Range("E96").Value = 801

# --- Re-ranked pairs: update data AND swap row order --------------------
# Paises Bajos (row 34) vs Emiratos Arabes Unidos (row 35):
# Emiratos' new total (50857) now exceeds Paises Bajos (50335), so
# Emiratos moves up to row 34 and Paises Bajos drops to row 35
# (Paises Bajos keeps its previous, unchanged figures).
$ws.Range("A34").Value = "Emiratos Arabes Unidos"
$ws.Range("B34").Value = 50857
$ws.Range("C34").Value = 716
$ws.Range("D34").Value = 39857
$ws.Range("E34").Value = 10679
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 3
$ws.Range("H34").Value = 321

$ws.Range("A35").Value = "Paises Bajos"
$ws.Range("B35").Value = 50335
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 0
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 6113

# Singapur (row 39) vs Oman (row 40):
# Oman's new total (45106) now exceeds Singapur (44664), so Oman moves
# up to row 39 and Singapur drops to row 40 (unchanged figures).
$ws.Range("A39").Value = "Oman"
$ws.Range("B39").Value = 45106
$ws.Range("C39").Value = 1177
$ws.Range("D39").Value = 26968
$ws.Range("E39").Value = 17935
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 10
$ws.Range("H39").Value = 203

$ws.Range("A40").Value = "Singapur"
$ws.Range("B40").Value = 44664
$ws.Range("C40").Value = 185
$ws.Range("D40").Value = 39769
$ws.Range("E40").Value = 4869
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 26

# Barein (row 49) vs Rumania (row 50):
# Rumania's new total (28582) now exceeds Barein (28410), so Rumania
# moves up to row 49 and Barein drops to row 50 (unchanged figures).
$ws.Range("A49").Value = "Rumania"
$ws.Range("B49").Value = 28582
$ws.Range("C49").Value = 416
$ws.Range("D49").Value = 19854
$ws.Range("E49").Value = 6997
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 23
$ws.Range("H49").Value = 1731

$ws.Range("A50").Value = "Barein"
$ws.Range("B50").Value = 28410
$ws.Range("C50").Value = 0
$ws.Range("D50").Value = 23318
$ws.Range("E50").Value = 4997
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 95

# Libano (row 116) vs Lituania (row 117):
# Lituania's new total (1831) now exceeds Libano (1830), so Lituania
# moves up to row 116 and Libano drops to row 117 (unchanged figures).
$ws.Range("A116").Value = "Lituania"
$ws.Range("B116").Value = 1831
$ws.Range("C116").Value = 3
$ws.Range("D116").Value = 1545
$ws.Range("E116").Value = 207
$ws.Range("F116").Value = 0
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 79

$ws.Range("A117").Value = "Libano"
$ws.Range("B117").Value = 1830
$ws.Range("C117").Value = 0
$ws.Range("D117").Value = 1292
$ws.Range("E117").Value = 503
$ws.Range("F117").Value = 0
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 35
